$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("E").ColumnWidth = 7.86
